$wb = $excel.ActiveWorkbook

# --- "Info" sheet: updated car model drag/lift coefficients + new low-drag configs ---
$info = $wb.Worksheets.Item("Info")

# Lift Coefficient CL: -2 -> -1.98
$info.Range("C8").Value = -1.98

# Drag Coefficient CD: -1.2 -> -1.33
$info.Range("C9").Value = -1.33

# Front Aero Distribution: 47 -> formula 100-56.3 (=43.7), a low-drag config split
$info.Range("C12").Formula = "=100-56.3"

# Frontal Area: 1.1 -> 1.15
$info.Range("C13").Value = 1.15

# Move the view/selection on the Info sheet to E12 and make it the active sheet/tab
$info.Activate() | Out-Null
$info.Range("E12").Select() | Out-Null

# --- "Torque Curve" sheet: selection stays at J48, but it's no longer the active tab ---
$torque = $wb.Worksheets.Item("Torque Curve")
$torque.Range("J48").Select() | Out-Null

# Re-activate Info so it ends up as the selected/active tab when saved
$info.Activate() | Out-Null
